# March 24 update 3
# Adds three new columns (M: renewd, N: PlanID, O: iteration) to Sheet1,
# filling every data row with the same constant values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1), matching the bold/centered/bordered header
# style used by the existing headers in B1:L1 (copy formatting from L1).
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Fill in the new columns for every data row (rows 2-19).
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 13).Value = "before"    # M: renewd
    $ws.Cells.Item($r, 14).Value = 20140060     # N: PlanID
    $ws.Cells.Item($r, 15).Value = 11           # O: iteration
}
